$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above row 59, pushing the existing rows 59-69 down to 61-71.
$ws.Range("A59:R60").EntireRow.Insert()

# Row 59 (new): Terminal Hortofruticola Agro Chillan - Perejil, Primera, 0,5-1 kg atado, Fecha 45077
$ws.Cells.Item(59, 1).Value = 7
$ws.Cells.Item(59, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(59, 3).Value = "Ñuble"
$ws.Cells.Item(59, 4).Value = 45077
$ws.Cells.Item(59, 5).Value = 16
$ws.Cells.Item(59, 6).Value = 100112044
$ws.Cells.Item(59, 7).Value = "Perejil"
$ws.Cells.Item(59, 8).Value = "Sin especificar"
$ws.Cells.Item(59, 9).Value = "Primera"
$ws.Cells.Item(59, 10).Value = 150
$ws.Cells.Item(59, 11).Value = 1000
$ws.Cells.Item(59, 12).Value = 1000
$ws.Cells.Item(59, 13).Value = 1000
$ws.Cells.Item(59, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(59, 15).Value = "Región del Maule"
$ws.Cells.Item(59, 16).Value = 1000
$ws.Cells.Item(59, 17).Value = 1
$ws.Cells.Item(59, 18).Value = "Hortaliza"

# Row 60 (new): Terminal Hortofruticola Agro Chillan - Perejil, Segunda, 0,5-1 kg atado, Fecha 45077
$ws.Cells.Item(60, 1).Value = 7
$ws.Cells.Item(60, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(60, 3).Value = "Ñuble"
$ws.Cells.Item(60, 4).Value = 45077
$ws.Cells.Item(60, 5).Value = 16
$ws.Cells.Item(60, 6).Value = 100112044
$ws.Cells.Item(60, 7).Value = "Perejil"
$ws.Cells.Item(60, 8).Value = "Sin especificar"
$ws.Cells.Item(60, 9).Value = "Segunda"
$ws.Cells.Item(60, 10).Value = 150
$ws.Cells.Item(60, 11).Value = 800
$ws.Cells.Item(60, 12).Value = 800
$ws.Cells.Item(60, 13).Value = 800
$ws.Cells.Item(60, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(60, 15).Value = "Región del Maule"
$ws.Cells.Item(60, 16).Value = 800
$ws.Cells.Item(60, 17).Value = 1
$ws.Cells.Item(60, 18).Value = "Hortaliza"
